$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '53.885.28'
$ws.Range('E2').Value = '  +0.16%  '

$ws.Range('D3').Value = '2.237.60'
$ws.Range('E3').Value = '  +0.00%  '

$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').Value = "'493.13"
$ws.Range('E5').Value = '  +1.15%  '

$ws.Range('D6').Value = "'126.79"
$ws.Range('E6').Value = '  +0.37%  '

$ws.Range('D7').Value = "'0.993"
$ws.Range('E7').Value = '  -0.67%  '

$ws.Range('E8').Value = '  +1.26%  '

$ws.Range('D9').Value = '2.276.72'
$ws.Range('E9').Value = '  +1.32%  '

$ws.Range('D10').Value = "'0.0946"
$ws.Range('E10').Value = '  +3.07%  '

$ws.Range('D11').Value = "'0.151"
$ws.Range('E11').Value = '  +2.01%  '

$ws.Range('E12').Value = '  +3.01%  '

$ws.Range('E13').Value = '  -2.50%  '

$ws.Range('D14').Value = '2.642.23'
$ws.Range('E14').Value = '  +0.20%  '

$ws.Range('D15').Value = "'21.72"
$ws.Range('E15').Value = '  +3.03%  '

$ws.Range('D16').Value = '53.849.94'
$ws.Range('E16').Value = '  +0.26%  '

$ws.Range('D17').Value = "'0.0000129"
$ws.Range('E17').Value = '  +0.87%  '

$ws.Range('D18').Value = '2.243.45'
$ws.Range('E18').Value = '  -0.28%  '

$ws.Range('D19').Value = "'10.04"
$ws.Range('E19').Value = '  +4.95%  '

$ws.Range('D20').Value = "'4.09"
$ws.Range('E20').Value = '  +3.05%  '

$ws.Range('D21').Value = "'6.46"
$ws.Range('E21').Value = '  +5.57%  '

$ws.Range('D22').Value = "'300.54"
$ws.Range('E22').Value = '  +0.49%  '

$ws.Range('E23').Value = '  -0.65%  '

$ws.Range('E24').Value = '  -2.40%  '

$ws.Range('D25').Value = "'62.25"
$ws.Range('E25').Value = '  -2.21%  '

$ws.Range('D26').Value = "'1.02"
$ws.Range('E26').Value = '  +1.95%  '

$ws.Range('E27').Value = '  +2.23%  '

$ws.Range('D28').Value = "'0.149"
$ws.Range('E28').Value = '  +4.08%  '

$ws.Range('D29').Value = '2.353.36'
$ws.Range('E29').Value = '  +0.32%  '

$ws.Range('D30').Value = "'7.06"
$ws.Range('E30').Value = '  +0.57%  '

$ws.Range('D31').Value = "'168.21"
$ws.Range('E31').Value = '  -0.62%  '

$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = "'1.59"
$ws.Range('E32').Value = '  +0.01%  '

$ws.Range('B33').Value = 'PEPE'
$ws.Range('C33').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D33').Value = '0.0₃0684'
$ws.Range('E33').Value = '  -0.55%  '

$ws.Range('D34').Value = "'5.85"
$ws.Range('E34').Value = '  +1.78%  '

$ws.Range('E35').Value = '  -0.08%  '

$ws.Range('D36').Value = "'0.991"
$ws.Range('E36').Value = '  -0.66%  '

$ws.Range('E37').Value = '  +1.09%  '

$ws.Range('D38').Value = "'17.57"
$ws.Range('E38').Value = '  +0.60%  '

$ws.Range('D39').Value = "'1.18"
$ws.Range('E39').Value = '  +2.52%  '

$ws.Range('D40').Value = "'0.861"
$ws.Range('E40').Value = '  +2.55%  '

$ws.Range('D41').Value = "'3.69"
$ws.Range('E41').Value = '  +3.40%  '

$ws.Range('D42').Value = "'35.34"
$ws.Range('E42').Value = '  -0.98%  '

$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = "'1.40"
$ws.Range('E43').Value = '  +2.86%  '

$ws.Range('B44').Value = 'PolygonEcosystemToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D44').Value = "'0.373"
$ws.Range('E44').Value = '  +2.32%  '

$ws.Range('E45').Value = '  +1.10%  '

$ws.Range('D46').Value = "'127.88"
$ws.Range('E46').Value = '  +4.29%  '

$ws.Range('D47').Value = "'4.90"
$ws.Range('E47').Value = '  +5.41%  '

$ws.Range('D48').Value = "'0.0886"
$ws.Range('E48').Value = '  +0.80%  '

$ws.Range('E49').Value = '  +0.88%  '

$ws.Range('D50').Value = "'237.65"
$ws.Range('E50').Value = '  +1.56%  '

$ws.Range('E51').Value = '  +2.53%  '
